$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 74, pushing the existing rows 74-178
# down to 75-179 (dimension grows from A1:R178 to A1:R179).
$ws.Rows("74:74").Insert()

# Populate the new row 74 with the new record's data. The columns that are
# constant for every data row in this sheet (A, B, C, E, F, G, H, I, R) are
# copied from the neighbouring row; the rest (D, J, K, L, M, N, O, P, Q) get
# the new values from the edit.
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 45100
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = 100112035
$ws.Range("G74").Value = "Bruselas (repollito)"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 35
$ws.Range("K74").Value = 28000
$ws.Range("L74").Value = 28000
$ws.Range("M74").Value = 28000
$ws.Range("N74").Value = "$/malla 15 kilos"
$ws.Range("O74").Value = "Región Metropolitana"
$ws.Range("P74").Value = 1867
$ws.Range("Q74").Value = 15
$ws.Range("R74").Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D74").NumberFormat = $ws.Range("D75").NumberFormat
